# Auto-generated Excel COM-interop script to apply diff changes
# to Sheets/Behemoth_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5081.8335
$ws.Range("I62").Value = 2163.8333
$ws.Range("J62").Value = 7999.8335
$ws.Range("K62").Value = 2163.8333
$ws.Range("L62").Value = 7999.8335
$ws.Range("M62").Value = -1539.8333
$ws.Range("N62").Value = -9247.833500000001
$ws.Range("H65").Value = 5081.8335
$ws.Range("I65").Value = 2163.8333
$ws.Range("J65").Value = 7999.8335
$ws.Range("K65").Value = 10819.1665
$ws.Range("L65").Value = 39999.1675
$ws.Range("M65").Value = -7699.166499999999
$ws.Range("N65").Value = -46239.1675
$ws.Range("H69").Value = 17818.834
$ws.Range("I69").Value = 7379.6
$ws.Range("J69").Value = 70015
$ws.Range("K69").Value = 22138.8
$ws.Range("L69").Value = 210045
$ws.Range("M69").Value = -21264.8
$ws.Range("N69").Value = -211793
$ws.Range("H72").Value = 17818.834
$ws.Range("I72").Value = 7379.6
$ws.Range("J72").Value = 70015
$ws.Range("K72").Value = 66416.40000000001
$ws.Range("L72").Value = 630135
$ws.Range("M72").Value = -62048.40000000001
$ws.Range("N72").Value = -638871
$ws.Range("H100").Value = 1020.8571
$ws.Range("I100").Value = 1020.8571
$ws.Range("K100").Value = 1020.8571
$ws.Range("M100").Value = -479.8570999999999
$ws.Range("H101").Value = 1793
$ws.Range("J101").Value = 1128.3334
$ws.Range("L101").Value = 3385.0002
$ws.Range("N101").Value = -6629.0002
$ws.Range("H131").Value = 2924.6667
$ws.Range("I131").Value = 1387
$ws.Range("K131").Value = 4161
$ws.Range("M131").Value = 879
$ws.Range("H137").Value = 3156.9846
$ws.Range("I137").Value = 2434.1628
$ws.Range("K137").Value = 7302.4884
$ws.Range("M137").Value = -4752.4884
$ws.Range("H138").Value = 1932.05
$ws.Range("J138").Value = 2447.9429
$ws.Range("L138").Value = 7343.8287
$ws.Range("N138").Value = -17623.8287

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12826301
$ws.Range("I32").Value = 15153918
$ws.Range("J32").Value = 24407.834
$ws.Range("K32").Value = 15153918
$ws.Range("L32").Value = 24407.834
$ws.Range("M32").Value = -15153631
$ws.Range("N32").Value = -24981.834
$ws.Range("H61").Value = 18755316
$ws.Range("I61").Value = 13518389
$ws.Range("K61").Value = 13518389
$ws.Range("M61").Value = -13518177
$ws.Range("H103").Value = 56798
$ws.Range("J103").Value = 56798
$ws.Range("L103").Value = 56798
$ws.Range("N103").Value = -59142
$ws.Range("H120").Value = 73990
$ws.Range("J120").Value = 73990
$ws.Range("L120").Value = 73990
$ws.Range("N120").Value = -83666
$ws.Range("H136").Value = 18755316
$ws.Range("I136").Value = 13518389
$ws.Range("K136").Value = 40555167
$ws.Range("M136").Value = -40552617

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3224.0833
$ws.Range("I22").Value = 2638.9
$ws.Range("J22").Value = 6150
$ws.Range("K22").Value = 2638.9
$ws.Range("L22").Value = 6150
$ws.Range("M22").Value = -2465.9
$ws.Range("N22").Value = -6496
$ws.Range("H94").Value = 1298.5385
$ws.Range("I94").Value = 733.25
$ws.Range("J94").Value = 1549.7778
$ws.Range("K94").Value = 733.25
$ws.Range("L94").Value = 1549.7778
$ws.Range("M94").Value = -282.25
$ws.Range("N94").Value = -2451.7778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 674091.75
$ws.Range("I62").Value = 916181.4399999999
$ws.Range("K62").Value = 916181.4399999999
$ws.Range("M62").Value = -915557.4399999999
$ws.Range("H65").Value = 674091.75
$ws.Range("I65").Value = 916181.4399999999
$ws.Range("K65").Value = 4580907.199999999
$ws.Range("M65").Value = -4577787.199999999
$ws.Range("H99").Value = 3142.7334
$ws.Range("I99").Value = 2638.5557
$ws.Range("J99").Value = 3899
$ws.Range("K99").Value = 2638.5557
$ws.Range("L99").Value = 3899
$ws.Range("M99").Value = -1140.5557
$ws.Range("N99").Value = -6895
$ws.Range("H126").Value = 3142.7334
$ws.Range("I126").Value = 2638.5557
$ws.Range("J126").Value = 3899
$ws.Range("K126").Value = 7915.6671
$ws.Range("L126").Value = 11697
$ws.Range("M126").Value = -5445.6671
$ws.Range("N126").Value = -16637
$ws.Range("H140").Value = 44100
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18208814
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224
$ws.Range("H7").Value = 1025.4
$ws.Range("I7").Value = 325.5
$ws.Range("J7").Value = 2075.25
$ws.Range("K7").Value = 976.5
$ws.Range("L7").Value = 6225.75
$ws.Range("M7").Value = -864.5
$ws.Range("N7").Value = -6449.75
$ws.Range("H23").Value = 1716.8334
$ws.Range("I23").Value = 1980
$ws.Range("J23").Value = 401
$ws.Range("K23").Value = 5940
$ws.Range("L23").Value = 1203
$ws.Range("M23").Value = -5705
$ws.Range("N23").Value = -1673
$ws.Range("H33").Value = 149
$ws.Range("I33").Value = 22.5
$ws.Range("J33").Value = 293.57144
$ws.Range("K33").Value = 135
$ws.Range("L33").Value = 1761.42864
$ws.Range("M33").Value = 148
$ws.Range("N33").Value = -2327.42864
$ws.Range("H34").Value = 62.8
$ws.Range("I34").Value = 62.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 188.4
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -104.4
$ws.Range("N34").ClearContents()
$ws.Range("H131").Value = 5653.2173
$ws.Range("J131").Value = 5653.2173
$ws.Range("L131").Value = 16959.6519
$ws.Range("N131").Value = -27039.6519
$ws.Range("H140").Value = 3151.739
$ws.Range("I140").Value = 2824.5
$ws.Range("K140").Value = 8473.5
$ws.Range("M140").Value = -3293.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 175.2
$ws.Range("I2").Value = 65.3
$ws.Range("K2").Value = 65.3
$ws.Range("M2").Value = 47.7
$ws.Range("H11").Value = 10827779
$ws.Range("I11").Value = 5070000
$ws.Range("K11").Value = 5070000
$ws.Range("M11").Value = -5069861
$ws.Range("H12").Value = 33926668
$ws.Range("I12").Value = 33926668
$ws.Range("K12").Value = 33926668
$ws.Range("M12").Value = -33926528
$ws.Range("H93").Value = 64251
$ws.Range("J93").Value = 64251
$ws.Range("L93").Value = 64251
$ws.Range("N93").Value = -67995
$ws.Range("H114").Value = 112326.336
$ws.Range("J114").Value = 112326.336
$ws.Range("L114").Value = 112326.336
$ws.Range("N114").Value = -121004.336
$ws.Range("H139").Value = 20296
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 151003.42
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 151003.42
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 151003.42
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -151227.42
$ws.Range("H40").Value = 4414.9165
$ws.Range("I40").Value = 2996.5
$ws.Range("J40").Value = 5833.3335
$ws.Range("K40").Value = 2996.5
$ws.Range("L40").Value = 5833.3335
$ws.Range("M40").Value = -2860.5
$ws.Range("N40").Value = -6105.3335
$ws.Range("H126").Value = 151003.42
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 151003.42
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 453010.26
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -457950.26
$ws.Range("H136").Value = 53776.93
$ws.Range("I136").Value = 8671.1875
$ws.Range("K136").Value = 26013.5625
$ws.Range("M136").Value = -23463.5625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 29000
$ws.Range("I37").Value = 29000
$ws.Range("K37").Value = 29000
$ws.Range("M37").Value = -28797
$ws.Range("H97").Value = 111000
$ws.Range("J97").Value = 111000
$ws.Range("L97").Value = 111000
$ws.Range("N97").Value = -112982
$ws.Range("H98").Value = 84773.8
$ws.Range("J98").Value = 84773.8
$ws.Range("L98").Value = 84773.8
$ws.Range("N98").Value = -90763.8
$ws.Range("H126").Value = 1665.7778
$ws.Range("I126").Value = 1249
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 3747
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -1277
$ws.Range("N126").Value = -19940
$ws.Range("H136").Value = 2698.5833
$ws.Range("I136").Value = 2418.1843
$ws.Range("K136").Value = 7254.5529
$ws.Range("M136").Value = -4704.5529
